# Updates cryptos list price/volume figures (cells D/E) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (these columns are pre-formatted,
# free-form strings in the source data -- e.g. "64.365.66" or "1.00" --
# not numbers) without leaving a stray NumberFormat/quote-prefix style
# behind on the cell.
function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '64.359.91'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.625.46'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws.Range("D5") '596.04'
$ws.Range("E5").Value = '  +0.05%  '
Set-TextValue $ws.Range("D6") '152.76'
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("E9").Value = '  +4.02%  '
Set-TextValue $ws.Range("D11") '0.395'
$ws.Range("E11").Value = '  +3.44%  '
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").Value = '3.096.34'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("E15").Value = '  +11.05%  '
$ws.Range("D16").Value = '64.261.82'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").Value = '2.624.77'
$ws.Range("E17").Value = '  +0.14%  '
Set-TextValue $ws.Range("D18") '12.29'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  +2.40%  '
Set-TextValue $ws.Range("D20") '349.15'
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("E21").Value = '  +2.91%  '
Set-TextValue $ws.Range("D22") '0.999'
$ws.Range("E22").Value = '  +0.09%  '
Set-TextValue $ws.Range("D23") '67.69'
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("E25").Value = '  -0.32%  '
Set-TextValue $ws.Range("D26") '1.66'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("E27").Value = '  +1.22%  '
Set-TextValue $ws.Range("D28") '548.33'
$ws.Range("E29").Value = '  +0.25%  '
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '0.0₃0911'
$ws.Range("E31").Value = '  +7.59%  '
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("E33").Value = '  +3.87%  '
Set-TextValue $ws.Range("D34") '5.51'
$ws.Range("E34").Value = '  +4.98%  '
Set-TextValue $ws.Range("D35") '6.22'
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").Value = '  +3.04%  '
Set-TextValue $ws.Range("D37") '165.48'
$ws.Range("E37").Value = '  -1.82%  '
Set-TextValue $ws.Range("D38") '20.11'
$ws.Range("E38").Value = '  +3.74%  '
Set-TextValue $ws.Range("D39") '2.01'
$ws.Range("E39").Value = '  +3.16%  '
Set-TextValue $ws.Range("D40") '0.999'
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -0.03%  '
Set-TextValue $ws.Range("D42") '168.51'
$ws.Range("E42").Value = '  +0.87%  '
Set-TextValue $ws.Range("D43") '41.86'
$ws.Range("E43").Value = '  +4.95%  '
Set-TextValue $ws.Range("D44") '4.12'
Set-TextValue $ws.Range("D45") '23.15'
$ws.Range("E45").Value = '  +7.40%  '
$ws.Range("E46").Value = '  +12.49%  '
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("E51").Value = '  -0.10%  '
